$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.748108439408092
$ws.Range("L2").Value = 0.670403870413533

$ws.Range("B3").Value = 0.538331317889367
$ws.Range("L3").Value = 0.62055734756295

$ws.Range("B4").Value = 0.458387576622686
$ws.Range("L4").Value = 0.703121504911223

$ws.Range("B5").Value = 0.340202418477655
$ws.Range("L5").Value = 0.66320354832362

$ws.Range("B6").Value = 0.33829822989737
$ws.Range("C6").Value = 0.366925945043565
$ws.Range("L6").Value = 0.407151516775038

$ws.Range("B7").Value = 0.329346428115222
$ws.Range("C7").Value = 0.523687810073323
$ws.Range("L7").Value = -0.0763319799713056

$ws.Range("B8").Value = 0.320249242216214
$ws.Range("L8").Value = 0.699570192636595

$ws.Range("B9").Value = 0.319190975705175
$ws.Range("L9").Value = 0.377306157800078

$ws.Range("B10").Value = 0.31547866481026
$ws.Range("L10").Value = -0.0300563026874272

$ws.Range("B11").Value = 0.010244201702362
$ws.Range("C11").Value = 0.0839217421719148
$ws.Range("L11").Value = -0.0145192493656787
